$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Refresh the "last saved" date shown by the datetimeFigureOut field that
#    lives in the Date Placeholder of the slide master and of every custom
#    layout (6/7/2018 -> 31/3/2019).
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.HasTextFrame) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Text -eq "6/7/2018") {
                    $tr.Text = "31/3/2019"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2. Rename "address book" -> "GradTrak" in the Undo/Redo activity diagram
#    on slide 1.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)

    # Shape ids are stable identifiers (p:cNvPr/@id) - use them instead of
    # Name, since the COM layer derives .Name for autoshapes from their
    # geometry preset rather than echoing the literal stored name.
    if ($shp.Id -eq 48) {
        # "TextBox 47" - "[command commits address book]"
        $origHeight = $shp.Height
        $shp.TextFrame.TextRange.Text = "[command commits GradTrak]"
        # Re-apply the original height: this textbox auto-fits to its
        # text, and the new wording reflows to (almost) the same box size;
        # pin it back to avoid a spurious sub-point size drift.
        $shp.Height = $origHeight
    }

    if ($shp.Id -eq 51) {
        # "Rectangle: Rounded Corners 50" - "Purge redundant states ..."
        $shp.TextFrame.TextRange.Text = "Purge redundant states and then save GradTrak to gradTrakStateList "
    }
}
